# Scheduled-runner market-data refresh: update the cached crafting-profit
# figures (currentAveragePrice / LevePrice / LeveProfit columns H-N) on the
# affected rows of each job sheet to the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2717.1606
$ws.Range("J17").Value = 2717.1606
$ws.Range("L17").Value = 8151.4818
$ws.Range("N17").Value = -8487.481800000001

$ws.Range("H80").Value = 19609962
$ws.Range("I80").Value = 30303910
$ws.Range("J80").Value = 4390.8335
$ws.Range("K80").Value = 90911730
$ws.Range("L80").Value = 13172.5005
$ws.Range("M80").Value = -90910732
$ws.Range("N80").Value = -15168.5005

$ws.Range("H83").Value = 19609962
$ws.Range("I83").Value = 30303910
$ws.Range("J83").Value = 4390.8335
$ws.Range("K83").Value = 272735190
$ws.Range("L83").Value = 39517.5015
$ws.Range("M83").Value = -272730198
$ws.Range("N83").Value = -49501.5015

$ws.Range("H88").Value = 561048.2
$ws.Range("I88").Value = 978747.1
$ws.Range("J88").Value = 4116.3335
$ws.Range("K88").Value = 978747.1
$ws.Range("L88").Value = 4116.3335
$ws.Range("M88").Value = -978341.1
$ws.Range("N88").Value = -4928.3335

$ws.Range("H91").Value = 561048.2
$ws.Range("I91").Value = 978747.1
$ws.Range("J91").Value = 4116.3335
$ws.Range("K91").Value = 978747.1
$ws.Range("L91").Value = 4116.3335
$ws.Range("M91").Value = -977343.1
$ws.Range("N91").Value = -6924.3335

$ws.Range("H98").Value = 131262.83
$ws.Range("I98").Value = 1960
$ws.Range("J98").Value = 777777
$ws.Range("K98").Value = 1960
$ws.Range("L98").Value = 777777
$ws.Range("M98").Value = -462
$ws.Range("N98").Value = -780773

$ws.Range("H122").Value = 131262.83
$ws.Range("I122").Value = 1960
$ws.Range("J122").Value = 777777
$ws.Range("K122").Value = 5880
$ws.Range("L122").Value = 2333331
$ws.Range("M122").Value = -3430
$ws.Range("N122").Value = -2338231

$ws.Range("H129").Value = 1322.7778
$ws.Range("J129").Value = 1175.069
$ws.Range("L129").Value = 3525.207
$ws.Range("N129").Value = -13525.207

$ws.Range("H137").Value = 4469.846
$ws.Range("I137").Value = 1265.875
$ws.Range("K137").Value = 3797.625
$ws.Range("M137").Value = -1247.625

$ws.Range("H138").Value = 2084.4421
$ws.Range("I138").Value = 1080.1041
$ws.Range("J138").Value = 3110.149
$ws.Range("K138").Value = 3240.3123
$ws.Range("L138").Value = 9330.447
$ws.Range("M138").Value = 1899.6877
$ws.Range("N138").Value = -19610.447

$ws.Range("H141").Value = 1468.6
$ws.Range("I141").Value = 744.3
$ws.Range("J141").Value = 7263
$ws.Range("K141").Value = 2232.9
$ws.Range("L141").Value = 21789
$ws.Range("M141").Value = 2947.1
$ws.Range("N141").Value = -32149

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9522.638999999999
$ws.Range("I32").Value = 8464.537
$ws.Range("K32").Value = 8464.537
$ws.Range("M32").Value = -8177.537

$ws.Range("H74").Value = 1153.4642
$ws.Range("I74").Value = 1050.589
$ws.Range("J74").Value = 1836.1818
$ws.Range("K74").Value = 1050.589
$ws.Range("L74").Value = 1836.1818
$ws.Range("M74").Value = -176.5889999999999
$ws.Range("N74").Value = -3584.1818

$ws.Range("H77").Value = 1153.4642
$ws.Range("I77").Value = 1050.589
$ws.Range("J77").Value = 1836.1818
$ws.Range("K77").Value = 5252.945
$ws.Range("L77").Value = 9180.909
$ws.Range("M77").Value = -884.9449999999997
$ws.Range("N77").Value = -17916.909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 414.14285
$ws.Range("I22").Value = 414.14285
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 414.14285
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -241.14285
$ws.Range("N22").ClearContents()

$ws.Range("H26").Value = 9837.888999999999
$ws.Range("I26").Value = 4756.8335
$ws.Range("K26").Value = 4756.8335
$ws.Range("M26").Value = -4464.8335

$ws.Range("H94").Value = 596.8182
$ws.Range("I94").Value = 579.2857
$ws.Range("J94").Value = 627.5
$ws.Range("K94").Value = 579.2857
$ws.Range("L94").Value = 627.5
$ws.Range("M94").Value = -128.2857
$ws.Range("N94").Value = -1529.5

$ws.Range("H107").Value = 2285.1365
$ws.Range("I107").Value = 2391.9092
$ws.Range("J107").Value = 2178.3635
$ws.Range("K107").Value = 2391.9092
$ws.Range("L107").Value = 2178.3635
$ws.Range("M107").Value = -471.9092000000001
$ws.Range("N107").Value = -6018.363499999999

$ws.Range("H134").Value = 3058.5076
$ws.Range("I134").Value = 1780.5555
$ws.Range("J134").Value = 3921.125
$ws.Range("K134").Value = 5341.666499999999
$ws.Range("L134").Value = 11763.375
$ws.Range("M134").Value = -2806.666499999999
$ws.Range("N134").Value = -16833.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4997.2876
$ws.Range("I31").Value = 3600.1428
$ws.Range("J31").Value = 5119.5376
$ws.Range("K31").Value = 3600.1428
$ws.Range("L31").Value = 5119.5376
$ws.Range("M31").Value = -3305.1428
$ws.Range("N31").Value = -5709.5376

$ws.Range("H34").Value = 4997.2876
$ws.Range("I34").Value = 3600.1428
$ws.Range("J34").Value = 5119.5376
$ws.Range("K34").Value = 3600.1428
$ws.Range("L34").Value = 5119.5376
$ws.Range("M34").Value = -3398.1428
$ws.Range("N34").Value = -5523.5376

$ws.Range("H58").Value = 1297.491
$ws.Range("I58").Value = 1025.8043
$ws.Range("J58").Value = 2686.111
$ws.Range("K58").Value = 1025.8043
$ws.Range("L58").Value = 2686.111
$ws.Range("M58").Value = -822.8043
$ws.Range("N58").Value = -3092.111

$ws.Range("H99").Value = 2772.5454
$ws.Range("I99").Value = 2398.2222
$ws.Range("J99").Value = 4457
$ws.Range("K99").Value = 2398.2222
$ws.Range("L99").Value = 4457
$ws.Range("M99").Value = -900.2222000000002
$ws.Range("N99").Value = -7453

$ws.Range("H126").Value = 2772.5454
$ws.Range("I126").Value = 2398.2222
$ws.Range("J126").Value = 4457
$ws.Range("K126").Value = 7194.6666
$ws.Range("L126").Value = 13371
$ws.Range("M126").Value = -4724.6666
$ws.Range("N126").Value = -18311

$ws.Range("H134").Value = 390340.4
$ws.Range("I134").Value = 898.6667
$ws.Range("J134").Value = 1558665.6
$ws.Range("K134").Value = 2696.0001
$ws.Range("L134").Value = 4675996.800000001
$ws.Range("M134").Value = -161.0001000000002
$ws.Range("N134").Value = -4681066.800000001

$ws.Range("H136").Value = 1297.491
$ws.Range("I136").Value = 1025.8043
$ws.Range("J136").Value = 2686.111
$ws.Range("K136").Value = 3077.4129
$ws.Range("L136").Value = 8058.333
$ws.Range("M136").Value = -527.4129000000003
$ws.Range("N136").Value = -13158.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4710.6665
$ws.Range("I113").Value = 7151.1333
$ws.Range("J113").Value = 643.2222
$ws.Range("K113").Value = 21453.3999
$ws.Range("L113").Value = 1929.6666
$ws.Range("M113").Value = -19283.3999
$ws.Range("N113").Value = -6269.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1353.7778
$ws.Range("I113").Value = 1541.4286
$ws.Range("J113").Value = 1234.3636
$ws.Range("K113").Value = 1541.4286
$ws.Range("L113").Value = 1234.3636
$ws.Range("M113").Value = 628.5714
$ws.Range("N113").Value = -5574.3636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2327.361
$ws.Range("I7").Value = 1901.5416
$ws.Range("J7").Value = 3179
$ws.Range("K7").Value = 1901.5416
$ws.Range("L7").Value = 3179
$ws.Range("M7").Value = -1789.5416
$ws.Range("N7").Value = -3403

$ws.Range("H93").Value = 1597.9
$ws.Range("I93").Value = 946.8
$ws.Range("K93").Value = 946.8
$ws.Range("M93").Value = 301.2

$ws.Range("H126").Value = 2327.361
$ws.Range("I126").Value = 1901.5416
$ws.Range("J126").Value = 3179
$ws.Range("K126").Value = 5704.6248
$ws.Range("L126").Value = 9537
$ws.Range("M126").Value = -3234.6248
$ws.Range("N126").Value = -14477

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1839990.8
$ws.Range("I126").Value = 2264219.8
$ws.Range("J126").Value = 1665
$ws.Range("K126").Value = 6792659.399999999
$ws.Range("L126").Value = 4995
$ws.Range("M126").Value = -6790189.399999999
$ws.Range("N126").Value = -9935

$ws.Range("H136").Value = 256798.9
$ws.Range("I136").Value = 344421.94
$ws.Range("J136").Value = 1895.5454
$ws.Range("K136").Value = 1033265.82
$ws.Range("L136").Value = 5686.6362
$ws.Range("M136").Value = -1030715.82
$ws.Range("N136").Value = -10786.6362

